# Updated the Chrome Option as we got failure for Origin related issue
$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("HomePage")

# C2 currently holds "new arrival" (test data value) -> correct to "new arrivals"
$ws.Range("C2").Value = "new arrivals"

# Move/set the active selection to C2 to match the saved workbook view state
$ws.Activate()
$ws.Range("C2").Select()
